$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete D2, add B2 and C2
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 43.137025930401123
$ws.Range("C2").Value = 21.834834647764225

# Row 3: delete B3, change C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 41.884430876923176

# Update selection to reflect new used range of interest
$ws.Range("B1:E3").Select()
